# daily auto push: 2026-01-20 09:44 UTC
# Inserts a new day's-first log row above the existing row 677, shifting the
# trailing rows (677-718) down by one. The row that used to fall off the end
# (old row 718, the second 2027/01/05 entry) is preserved as the new last
# row (719), and the freshly inserted row 677 carries the new 2026/01/20
# entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 677:718 down to 678:719, duplicating row 677's formatting into
# the newly opened row (there is no cell styling on this range, so this is a
# plain shift).
$ws.Rows("677:677").Insert()

# New row 677: 2026/01/20 is a Tuesday (same weekday glyph as 2026/12/29 in
# this log - both "火").
$ws.Cells.Item(677, 1).Value = "2026/01/20"
$ws.Cells.Item(677, 2).Value = "火"
$ws.Cells.Item(677, 3).Value = 16
$ws.Cells.Item(677, 4).Value = 201
